$d = $word.ActiveDocument

# First paragraph: the one holding the "**ID__AFFARS_MP5301_6__ID**" placeholder.
$p = $d.Paragraphs(1)

# Remove the trailing standalone-space run by replacing the paragraph's
# text (placeholder + trailing space) with just the placeholder text.
$p.Range.Find.Execute("**ID__AFFARS_MP5301_6__ID** ", $true, $false, $false, $false, $false, `
                       $true, 1, $false, "**ID__AFFARS_MP5301_6__ID**", 2) | Out-Null

# Add a paragraph border (top/left/bottom/right) with 5pt text-distance,
# matching <w:pBdr><w:top w:space="5"/>...</w:pBdr>.
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5

# Widen the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p.Format.LeftIndent = 11.25
